# Repeater method edit: append the 4th "bio stock" data row (row 32) and
# backfill the PriceChange/UpDown comparison (X31/Y31) for the row that
# used to be last (row 31), now that row 32 exists to compare against.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Backfill row 31: PriceChange (X31) / UpDown (Y31) ----
$ws.Cells.Item(31, 24).Value = -0.21000000000000085   # X31 PriceChange
$ws.Cells.Item(31, 25).Value = "Down"                  # Y31 UpDown

# ---- New row 32 ----
$ws.Cells.Item(32, 1).Value = 42651.425162037034       # A32 Date
$ws.Cells.Item(32, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(32, 2).Value = 13                       # B32 ScoreFinal
$ws.Cells.Item(32, 3).Value = "Buy"                    # C32 Verdict
$ws.Cells.Item(32, 4).Value = 56                       # D32 totalSentiment
$ws.Cells.Item(32, 5).Value = 1623                     # E32 wordCount
$ws.Cells.Item(32, 6).Value = 181                      # F32 sentenceCount
$ws.Cells.Item(32, 7).Value = 69                       # G32 posWordPercentage
$ws.Cells.Item(32, 8).Value = 30                       # H32 negWordPercentage
$ws.Cells.Item(32, 9).Value = 94                       # I32 posPhrasePercentage
$ws.Cells.Item(32, 10).Value = 5                       # J32 negPhrasePercentage
$ws.Cells.Item(32, 11).Value = 4610                    # K32 ElapsedMs
$ws.Cells.Item(32, 12).Value = 23                      # L32 posWordCount
$ws.Cells.Item(32, 13).Value = 10                      # M32 negWordCount
$ws.Cells.Item(32, 14).Value = 16                      # N32 positivePhraseCount
$ws.Cells.Item(32, 15).Value = 1                       # O32 negativePhraseCount
$ws.Cells.Item(32, 16).Value = "Bag"                   # P32 Method
$ws.Cells.Item(32, 17).Value = 47.963765586266284      # Q32 RSI
$ws.Cells.Item(32, 18).Value = 0.49                    # R32 PEG

$ws.Cells.Item(32, 19).Value = 0.0521                  # S32 200Moving%
$ws.Cells.Item(32, 19).NumberFormat = "0.00%"

$ws.Cells.Item(32, 20).Value = -0.0214                 # T32 50Moving%
$ws.Cells.Item(32, 20).NumberFormat = "0.00%"

$ws.Cells.Item(32, 21).Value = 2.28                    # U32 PriceBook
$ws.Cells.Item(32, 22).Value = "N/A"                   # V32 Dividend
$ws.Cells.Item(32, 23).Value = 0                       # W32 Bollinger
# (no X32/Y32 -- PriceChange/UpDown only computed once a following row exists)

# ---- Re-apply the best-fit column widths the repeater recalculates on
#      every write (column A grows because the sheet has the longest date
#      column; B..M grow/shrink to fit the new row's values). ----
$ws.Columns.Item(1).ColumnWidth = 14.0234375
$ws.Columns.Item(2).ColumnWidth = 9.3203125
$ws.Columns.Item(3).ColumnWidth = 6.87890625
$ws.Columns.Item(4).ColumnWidth = 13.73828125
$ws.Columns.Item(5).ColumnWidth = 10.0234375
$ws.Columns.Item(6).ColumnWidth = 13.73828125
$ws.Range("G1:H1").ColumnWidth = 15.59375
$ws.Range("I1:J1").ColumnWidth = 16.59375
$ws.Range("L1:M1").ColumnWidth = 11.59375
